$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C ("Förändrad") from 45184 to 45186 for all data rows (2..305)
for ($r = 2; $r -le 305; $r++) {
    $c = $ws.Cells.Item($r, 3)
    if ($c.Value2 -eq 45184) {
        $c.Value = 45186
    }
}

# 2) Add the friendly-name second argument to the HYPERLINK formulas in columns
#    S..Y (19..25) for the rows that have them (rows 2..13).
$hyperlinkCols = @(19, 20, 21, 22, 23, 24, 25)
for ($r = 2; $r -le 13; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f.EndsWith('")')) {
                $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $name + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}

# 3) Row 305 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(305).RowHeight = 15

# 4) Append three new rows (306, 307, 308) with data
$newRows = @(
    @{ Row = 306; A = "A 43546-2023"; B = 45184; C = 45186; G = 1.5 },
    @{ Row = 307; A = "A 43495-2023"; B = 45184; C = 45186; G = 2.7 },
    @{ Row = 308; A = "A 43611-2023"; B = 45184; C = 45186; G = 1.2 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    $ws.Cells.Item($r, 1).Value = $rowData.A

    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value = "ÖSTERGÖTLANDS LÄN"
    $ws.Cells.Item($r, 5).Value = "SÖDERKÖPING"

    $ws.Cells.Item($r, 7).Value = $rowData.G
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = 0
    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 0
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
    $ws.Cells.Item($r, 14).Value = 0
    $ws.Cells.Item($r, 15).Value = 0
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0

    $ws.Cells.Item($r, 18).Value = ""
    $ws.Cells.Item($r, 18).WrapText = $true
}

# Rows 306 and 307 have an explicit row height, like most other rows; row 308
# does not (matching the source data).
$ws.Rows.Item(306).RowHeight = 15
$ws.Rows.Item(307).RowHeight = 15

Write-Output "done"
